$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added at the top of the data (row 73), pushing the
# existing records (old rows 73-155) down by one (new rows 74-156).
$ws.Rows("73:73").Insert()

# Populate the newly inserted row 73 with the new record's values.
$ws.Range("A73").Value = 10
$ws.Range("B73").Value = 'Vega Modelo de Temuco'
$ws.Range("C73").Value = 'La Araucanía'
$ws.Range("D73").Value = 44803
$ws.Range("E73").Value = 9
$ws.Range("F73").Value = 100112031
$ws.Range("G73").Value = 'Poroto verde'
$ws.Range("H73").Value = 'Sin especificar'
$ws.Range("I73").Value = 'Primera'
$ws.Range("J73").Value = 80
$ws.Range("K73").Value = 35000
$ws.Range("L73").Value = 35000
$ws.Range("M73").Value = 35000
$ws.Range("N73").Value = '$/malla 25 kilos'
$ws.Range("O73").Value = 'Provincia de Limarí'
$ws.Range("P73").Value = 1400
$ws.Range("Q73").Value = 25
$ws.Range("R73").Value = 'Hortaliza'
